$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes -------------------------------------------------
# The first addressing-mode table (row 5-6) changes its name from
# "Indirect" to "Direct", and since the direct mode only uses one
# register, the R2 field becomes unused.
$ws.Range("A5").Value = "Direct"
$ws.Range("D5").Value = "NOT USED"

# --- Formatting changes ------------------------------------------------
# The title cell (B2) picks up the same shaded / top+bottom-bordered,
# centered, non-bold look already used by the grid filler cells (e.g.
# F12), instead of its previous plain bold style.
$ws.Range("F12").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection change ----------------------------------------------------
$ws.Range("J15").Select() | Out-Null
